$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric strings (e.g. "221.91", "1.00") are not
# auto-converted to floating point numbers by Excel - force column D
# (and E, for safety) to Text format before writing values.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "33.957.51"
$ws.Range("E2").Value = "  -2.09%  "

$ws.Range("D3").Value = "1.790.33"
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "221.91"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").Value = "0.555"
$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("D8").Value = "32.31"
$ws.Range("E8").Value = "  -0.41%  "

$ws.Range("D9").Value = "0.283"
$ws.Range("E9").Value = "  +0.33%  "

$ws.Range("D10").Value = "0.0712"
$ws.Range("E10").Value = "  +2.71%  "

$ws.Range("E11").Value = "  -0.61%  "

$ws.Range("D12").Value = "2.049.37"
$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").Value = "1.794.38"
$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").Value = "10.77"
$ws.Range("E14").Value = "  -1.84%  "

$ws.Range("D15").Value = "0.624"
$ws.Range("E15").Value = "  -1.71%  "

$ws.Range("D16").Value = "33.945.66"
$ws.Range("E16").Value = "  -2.03%  "

$ws.Range("D17").Value = "4.17"
$ws.Range("E17").Value = "  -3.29%  "

$ws.Range("D18").Value = "67.94"
$ws.Range("E18").Value = "  -1.64%  "

$ws.Range("D19").Value = "243.90"
$ws.Range("E19").Value = "  -4.31%  "

$ws.Range("D20").Value = "0.0₃0782"
$ws.Range("E20").Value = "  -0.76%  "

$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").Value = "10.76"
$ws.Range("E22").Value = "  +1.43%  "

$ws.Range("D23").Value = "4.08"
$ws.Range("E23").Value = "  -2.71%  "

$ws.Range("D24").Value = "2.11"
$ws.Range("E24").Value = "  -1.64%  "

$ws.Range("D25").Value = "158.11"
$ws.Range("E25").Value = "  -1.41%  "

$ws.Range("D26").Value = "16.33"
$ws.Range("E26").Value = "  -0.41%  "

$ws.Range("D27").Value = "7.03"
$ws.Range("E27").Value = "  -1.68%  "

$ws.Range("E28").Value = "  -1.70%  "

$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").Value = "0.0520"
$ws.Range("E30").Value = "  -2.20%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.20"
$ws.Range("E31").Value = "  +0.66%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "3.67"
$ws.Range("E32").Value = "  -3.54%  "

$ws.Range("D33").Value = "3.48"
$ws.Range("E33").Value = "  -3.19%  "

$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  -2.84%  "

$ws.Range("D35").Value = "1.395.51"
$ws.Range("E35").Value = "  -3.16%  "

$ws.Range("D36").Value = "0.638"
$ws.Range("E36").Value = "  -0.35%  "

$ws.Range("D37").Value = "1.05"
$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("D38").Value = "0.0185"
$ws.Range("E38").Value = "  -3.61%  "

$ws.Range("D39").Value = "0.924"
$ws.Range("E39").Value = "  -0.37%  "

$ws.Range("D40").Value = "79.29"
$ws.Range("E40").Value = "  -7.50%  "

$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "2.34"
$ws.Range("E41").Value = "  +0.57%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "2.71"
$ws.Range("E42").Value = "  -2.69%  "

$ws.Range("D43").Value = "2.12"
$ws.Range("E43").Value = "  +1.51%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "0.0496"
$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "1.05"
$ws.Range("E45").Value = "  -1.22%  "

$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "106.96"
$ws.Range("E46").Value = "  +1.22%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "5.86"
$ws.Range("E47").Value = "  -2.34%  "

$ws.Range("D48").Value = "1.948.39"
$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "11.95"
$ws.Range("E50").Value = "  -1.80%  "

$ws.Range("D51").Value = "0.0₆0126"
$ws.Range("E51").Value = "  +4.08%  "
